$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(6, 2).Value = "[]"
$ws.Cells.Item(6, 3).Value = "[]"
$ws.Cells.Item(6, 4).Value = "[]"
$ws.Cells.Item(6, 5).Value = "[]"
$ws.Cells.Item(6, 6).Value = 0

$ws.Cells.Item(7, 2).Value = "[]"
$ws.Cells.Item(7, 3).Value = "[]"
$ws.Cells.Item(7, 4).Value = "[]"
$ws.Cells.Item(7, 5).Value = "[]"
$ws.Cells.Item(7, 6).Value = 0

$ws.Cells.Item(8, 2).Value = "['Fiat 500 E', 'Others']"
$ws.Cells.Item(8, 3).Value = "[0.4, 0.15]"
$ws.Cells.Item(8, 4).Value = "[0.7500000000000001, 0.8500000000000002]"
$ws.Cells.Item(8, 5).Value = "[8.400000000000002, 36.84138888888889]"
$ws.Cells.Item(8, 6).Value = 45.24138888888889

$ws.Cells.Item(9, 2).Value = "['Others', 'Audi E-TRON', 'Hyundai KONA 39 kWh']"
$ws.Cells.Item(9, 3).Value = "[0.25, 0.25, 0.4999999999999999]"
$ws.Cells.Item(9, 4).Value = "[1.0, 0.9500000000000003, 0.8000000000000002]"
$ws.Cells.Item(9, 5).Value = "[39.47291666666666, 59.50000000000002, 11.70000000000001]"
$ws.Cells.Item(9, 6).Value = 110.6729166666667

$ws.Cells.Item(12, 2).Value = "[]"
$ws.Cells.Item(12, 3).Value = "[]"
$ws.Cells.Item(12, 4).Value = "[]"
$ws.Cells.Item(12, 5).Value = "[]"
$ws.Cells.Item(12, 6).Value = 0

$ws.Cells.Item(13, 2).Value = "['VW ID.3', 'VW E-UP', 'Others', 'Polestar 2', 'SKODA ENYAQ 58kWh']"
$ws.Cells.Item(13, 3).Value = "[0.2, 0.3, 0.1, 0.25, 0.35]"
$ws.Cells.Item(13, 4).Value = "[0.7000000000000001, 0.8500000000000002, 0.8500000000000002, 0.9500000000000003, 0.8500000000000002]"
$ws.Cells.Item(13, 5).Value = "[29.0, 20.24000000000001, 39.47291666666668, 52.50000000000002, 29.000000000000014]"
$ws.Cells.Item(13, 6).Value = 170.2129166666667

$ws.Cells.Item(16, 2).Value = "['VW ID.3', 'VW ID.5', 'Others', 'Opel MOKKA']"
$ws.Cells.Item(16, 3).Value = "[0.05, 0.3, 0.1, 0.25]"
$ws.Cells.Item(16, 4).Value = "[0.8500000000000002, 0.65, 1.0, 0.7000000000000001]"
$ws.Cells.Item(16, 5).Value = "[46.400000000000006, 26.950000000000003, 47.3675, 20.250000000000004]"
$ws.Cells.Item(16, 6).Value = 140.9675

$ws.Cells.Item(17, 2).Value = "['Others']"
$ws.Cells.Item(17, 3).Value = "[0.2]"
$ws.Cells.Item(17, 4).Value = "[0.9500000000000003]"
$ws.Cells.Item(17, 5).Value = "[39.47291666666668]"
$ws.Cells.Item(17, 6).Value = 39.47291666666668

$ws.Cells.Item(18, 2).Value = "[]"
$ws.Cells.Item(18, 3).Value = "[]"
$ws.Cells.Item(18, 4).Value = "[]"
$ws.Cells.Item(18, 5).Value = "[]"
$ws.Cells.Item(18, 6).Value = 0

$ws.Cells.Item(19, 2).Value = "[]"
$ws.Cells.Item(19, 3).Value = "[]"
$ws.Cells.Item(19, 4).Value = "[]"
$ws.Cells.Item(19, 5).Value = "[]"
$ws.Cells.Item(19, 6).Value = 0

$ws.Cells.Item(30, 2).Value = "[]"
$ws.Cells.Item(30, 3).Value = "[]"
$ws.Cells.Item(30, 4).Value = "[]"
$ws.Cells.Item(30, 5).Value = "[]"
$ws.Cells.Item(30, 6).Value = 0

$ws.Cells.Item(31, 2).Value = "[]"
$ws.Cells.Item(31, 3).Value = "[]"
$ws.Cells.Item(31, 4).Value = "[]"
$ws.Cells.Item(31, 5).Value = "[]"
$ws.Cells.Item(31, 6).Value = 0

$ws.Cells.Item(32, 2).Value = "['VW ID.4', 'Others']"
$ws.Cells.Item(32, 3).Value = "[0.2, 0.25]"
$ws.Cells.Item(32, 4).Value = "[1.0, 0.8000000000000002]"
$ws.Cells.Item(32, 5).Value = "[61.6, 28.946805555555564]"
$ws.Cells.Item(32, 6).Value = 90.54680555555557

$ws.Cells.Item(33, 2).Value = "['Hyundai KONA 64 kWh', 'Hyundai IONIQ5 77kWh']"
$ws.Cells.Item(33, 3).Value = "[0.25, 0.3]"
$ws.Cells.Item(33, 4).Value = "[0.7500000000000001, 0.8500000000000002]"
$ws.Cells.Item(33, 5).Value = "[32.00000000000001, 42.35000000000002]"
$ws.Cells.Item(33, 6).Value = 74.35000000000002

$ws.Cells.Item(34, 2).Value = "['Others']"
$ws.Cells.Item(34, 3).Value = "[0.4]"
$ws.Cells.Item(34, 4).Value = "[0.65]"
$ws.Cells.Item(34, 5).Value = "[13.157638888888888]"
$ws.Cells.Item(34, 6).Value = 13.15763888888889

$ws.Cells.Item(36, 2).Value = "[]"
$ws.Cells.Item(36, 3).Value = "[]"
$ws.Cells.Item(36, 4).Value = "[]"
$ws.Cells.Item(36, 5).Value = "[]"
$ws.Cells.Item(36, 6).Value = 0

$ws.Cells.Item(37, 2).Value = "['Opel CORSA', 'Hyundai IONIQ5 58kWh', 'Others', 'Renault ZOE', 'Audi E-TRON']"
$ws.Cells.Item(37, 3).Value = "[0.35, 0.2, 0.2, 0.1, 0.35]"
$ws.Cells.Item(37, 4).Value = "[0.7000000000000001, 0.7000000000000001, 0.6, 0.9500000000000003, 0.8500000000000002]"
$ws.Cells.Item(37, 5).Value = "[15.750000000000004, 29.0, 21.05222222222222, 44.20000000000002, 42.50000000000002]"
$ws.Cells.Item(37, 6).Value = 152.5022222222223

$ws.Cells.Item(38, 2).Value = "[]"
$ws.Cells.Item(38, 3).Value = "[]"
$ws.Cells.Item(38, 4).Value = "[]"
$ws.Cells.Item(38, 5).Value = "[]"
$ws.Cells.Item(38, 6).Value = 0

$ws.Cells.Item(40, 2).Value = "['Renault ZOE', 'Smart FORTWO']"
$ws.Cells.Item(40, 3).Value = "[0.35, 0.05]"
$ws.Cells.Item(40, 4).Value = "[1.0, 0.65]"
$ws.Cells.Item(40, 5).Value = "[33.800000000000004, 10.56]"
$ws.Cells.Item(40, 6).Value = 44.36000000000001

$ws.Cells.Item(41, 2).Value = "['Fiat 500 E', 'Smart FORTWO', 'Audi E-TRON']"
$ws.Cells.Item(41, 3).Value = "[0.35, 0.4, 0.1]"
$ws.Cells.Item(41, 4).Value = "[0.7500000000000001, 0.9500000000000003, 0.8500000000000002]"
$ws.Cells.Item(41, 5).Value = "[9.600000000000003, 9.680000000000005, 63.75000000000002]"
$ws.Cells.Item(41, 6).Value = 83.03000000000003

$ws.Cells.Item(42, 2).Value = "[]"
$ws.Cells.Item(42, 3).Value = "[]"
$ws.Cells.Item(42, 4).Value = "[]"
$ws.Cells.Item(42, 5).Value = "[]"
$ws.Cells.Item(42, 6).Value = 0

$ws.Cells.Item(43, 2).Value = "[]"
$ws.Cells.Item(43, 3).Value = "[]"
$ws.Cells.Item(43, 4).Value = "[]"
$ws.Cells.Item(43, 5).Value = "[]"
$ws.Cells.Item(43, 6).Value = 0

$ws.Cells.Item(54, 2).Value = "[]"
$ws.Cells.Item(54, 3).Value = "[]"
$ws.Cells.Item(54, 4).Value = "[]"
$ws.Cells.Item(54, 5).Value = "[]"
$ws.Cells.Item(54, 6).Value = 0

$ws.Cells.Item(55, 2).Value = "['VW ID.5']"
$ws.Cells.Item(55, 3).Value = "[0.25]"
$ws.Cells.Item(55, 4).Value = "[0.9000000000000002]"
$ws.Cells.Item(55, 5).Value = "[50.05000000000002]"
$ws.Cells.Item(55, 6).Value = 50.05000000000002

$ws.Cells.Item(56, 2).Value = "['VW ID.5', 'Opel MOKKA']"
$ws.Cells.Item(56, 3).Value = "[0.3, 0.35]"
$ws.Cells.Item(56, 4).Value = "[0.7000000000000001, 1.0]"
$ws.Cells.Item(56, 5).Value = "[30.800000000000004, 29.25]"
$ws.Cells.Item(56, 6).Value = 60.05

$ws.Cells.Item(57, 2).Value = "['Opel MOKKA', 'Renault ZOE']"
$ws.Cells.Item(57, 3).Value = "[0.35, 0.2]"
$ws.Cells.Item(57, 4).Value = "[0.7000000000000001, 0.65]"
$ws.Cells.Item(57, 5).Value = "[15.750000000000004, 23.400000000000002]"
$ws.Cells.Item(57, 6).Value = 39.15000000000001

$ws.Cells.Item(60, 2).Value = "['Others', 'Others', 'Tesla MODEL 3', 'VW E-UP', 'Others']"
$ws.Cells.Item(60, 3).Value = "[0.1, 0.4, 0.2, 0.35, 0.1]"
$ws.Cells.Item(60, 4).Value = "[0.8000000000000002, 0.7500000000000001, 0.8500000000000002, 0.7000000000000001, 0.8000000000000002]"
$ws.Cells.Item(60, 5).Value = "[36.84138888888889, 18.420694444444447, 32.50000000000001, 12.880000000000003, 36.84138888888889]"
$ws.Cells.Item(60, 6).Value = 137.4834722222222

$ws.Cells.Item(61, 2).Value = "[]"
$ws.Cells.Item(61, 3).Value = "[]"
$ws.Cells.Item(61, 4).Value = "[]"
$ws.Cells.Item(61, 5).Value = "[]"
$ws.Cells.Item(61, 6).Value = 0

$ws.Cells.Item(62, 2).Value = "[]"
$ws.Cells.Item(62, 3).Value = "[]"
$ws.Cells.Item(62, 4).Value = "[]"
$ws.Cells.Item(62, 5).Value = "[]"
$ws.Cells.Item(62, 6).Value = 0

$ws.Cells.Item(64, 2).Value = "['Tesla MODEL 3', 'Others', 'Tesla MODEL 3']"
$ws.Cells.Item(64, 3).Value = "[0.2, 0.15, 0.3]"
$ws.Cells.Item(64, 4).Value = "[0.8000000000000002, 0.7000000000000001, 0.8500000000000002]"
$ws.Cells.Item(64, 5).Value = "[30.000000000000004, 28.946805555555557, 27.500000000000014]"
$ws.Cells.Item(64, 6).Value = 86.44680555555557

$ws.Cells.Item(65, 2).Value = "['TESLA MODEL Y', 'Others']"
$ws.Cells.Item(65, 3).Value = "[0.4, 0.2]"
$ws.Cells.Item(65, 4).Value = "[0.65, 0.8000000000000002]"
$ws.Cells.Item(65, 5).Value = "[18.75, 31.578333333333337]"
$ws.Cells.Item(65, 6).Value = 50.32833333333333

$ws.Cells.Item(66, 2).Value = "[]"
$ws.Cells.Item(66, 3).Value = "[]"
$ws.Cells.Item(66, 4).Value = "[]"
$ws.Cells.Item(66, 5).Value = "[]"
$ws.Cells.Item(66, 6).Value = 0

$ws.Cells.Item(67, 2).Value = "[]"
$ws.Cells.Item(67, 3).Value = "[]"
$ws.Cells.Item(67, 4).Value = "[]"
$ws.Cells.Item(67, 5).Value = "[]"
$ws.Cells.Item(67, 6).Value = 0

$ws.Cells.Item(78, 2).Value = "[]"
$ws.Cells.Item(78, 3).Value = "[]"
$ws.Cells.Item(78, 4).Value = "[]"
$ws.Cells.Item(78, 5).Value = "[]"
$ws.Cells.Item(78, 6).Value = 0

$ws.Cells.Item(79, 2).Value = "['Others']"
$ws.Cells.Item(79, 3).Value = "[0.1]"
$ws.Cells.Item(79, 4).Value = "[0.9000000000000002]"
$ws.Cells.Item(79, 5).Value = "[42.104444444444454]"
$ws.Cells.Item(79, 6).Value = 42.10444444444445

$ws.Cells.Item(80, 2).Value = "['Others']"
$ws.Cells.Item(80, 3).Value = "[0.1]"
$ws.Cells.Item(80, 4).Value = "[0.8000000000000002]"
$ws.Cells.Item(80, 5).Value = "[36.84138888888889]"
$ws.Cells.Item(80, 6).Value = 36.84138888888889

$ws.Cells.Item(81, 2).Value = "['MINI Cooper SE', 'TESLA MODEL Y']"
$ws.Cells.Item(81, 3).Value = "[0.3, 0.2]"
$ws.Cells.Item(81, 4).Value = "[0.8000000000000002, 0.8000000000000002]"
$ws.Cells.Item(81, 5).Value = "[14.450000000000006, 45.00000000000001]"
$ws.Cells.Item(81, 6).Value = 59.45000000000002

$ws.Cells.Item(82, 2).Value = "['VW ID.3']"
$ws.Cells.Item(82, 3).Value = "[0.1]"
$ws.Cells.Item(82, 4).Value = "[0.8000000000000002]"
$ws.Cells.Item(82, 5).Value = "[40.60000000000001]"
$ws.Cells.Item(82, 6).Value = 40.60000000000001

$ws.Cells.Item(83, 2).Value = "[]"
$ws.Cells.Item(83, 3).Value = "[]"
$ws.Cells.Item(83, 4).Value = "[]"
$ws.Cells.Item(83, 5).Value = "[]"
$ws.Cells.Item(83, 6).Value = 0

$ws.Cells.Item(84, 2).Value = "['Opel CORSA']"
$ws.Cells.Item(84, 3).Value = "[0.35]"
$ws.Cells.Item(84, 4).Value = "[0.8500000000000002]"
$ws.Cells.Item(84, 5).Value = "[22.50000000000001]"
$ws.Cells.Item(84, 6).Value = 22.50000000000001

$ws.Cells.Item(85, 2).Value = "['Audi Q4', 'Others', 'Hyundai KONA 39 kWh', 'Fiat 500 E']"
$ws.Cells.Item(85, 3).Value = "[0.4, 0.1, 0.2, 0.35]"
$ws.Cells.Item(85, 4).Value = "[0.7000000000000001, 0.8500000000000002, 0.8000000000000002, 0.9500000000000003]"
$ws.Cells.Item(85, 5).Value = "[22.98, 39.47291666666668, 23.400000000000002, 14.400000000000007]"
$ws.Cells.Item(85, 6).Value = 100.2529166666667

$ws.Cells.Item(88, 2).Value = "['Others', 'Others']"
$ws.Cells.Item(88, 3).Value = "[0.1, 0.3]"
$ws.Cells.Item(88, 4).Value = "[0.8500000000000002, 0.9500000000000003]"
$ws.Cells.Item(88, 5).Value = "[39.47291666666668, 34.20986111111113]"
$ws.Cells.Item(88, 6).Value = 73.6827777777778

$ws.Cells.Item(89, 2).Value = "['Renault ZOE', 'Fiat 500 E', 'Smart FORTWO']"
$ws.Cells.Item(89, 3).Value = "[0.25, 0.3, 0.3]"
$ws.Cells.Item(89, 4).Value = "[0.7000000000000001, 0.8500000000000002, 0.9000000000000002]"
$ws.Cells.Item(89, 5).Value = "[23.400000000000002, 13.200000000000006, 10.560000000000006]"
$ws.Cells.Item(89, 6).Value = 47.16000000000001

$ws.Cells.Item(90, 2).Value = "['Renault ZOE']"
$ws.Cells.Item(90, 3).Value = "[0.25]"
$ws.Cells.Item(90, 4).Value = "[0.7000000000000001]"
$ws.Cells.Item(90, 5).Value = "[23.400000000000002]"
$ws.Cells.Item(90, 6).Value = 23.4

$ws.Cells.Item(91, 2).Value = "[]"
$ws.Cells.Item(91, 3).Value = "[]"
$ws.Cells.Item(91, 4).Value = "[]"
$ws.Cells.Item(91, 5).Value = "[]"
$ws.Cells.Item(91, 6).Value = 0

$ws.Cells.Item(101, 2).Value = "[]"
$ws.Cells.Item(101, 3).Value = "[]"
$ws.Cells.Item(101, 4).Value = "[]"
$ws.Cells.Item(101, 5).Value = "[]"
$ws.Cells.Item(101, 6).Value = 0
